$wb = $excel.ActiveWorkbook

# --- Rename the "Include from ActMood" sheet to "Include #0" ---
$wsInclude = $wb.Worksheets.Item("Include from ActMood")
$wsInclude.Name = "Include #0"

# --- Update the Metadata sheet ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Update Version value (B3)
$ws1.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (B8)
$ws1.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row after "Contact" (row 10) for the "Jurisdiction" property,
# pushing Description/Purpose/Copyright/Immutable down by one row.
$ws1.Rows("11:11").Insert()

# Copy the formatting from the row below (now row 12, formerly row 11) so the
# new row matches the other data rows' style.
$ws1.Range("A12:B12").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

Write-Output "done"
